$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (rows 2-5), matching the "Numero / Nombre Completo / Local / Quien Autoriza / Fecha de Alta" header
$data = @(
    @("35115887", "Analia Belen Miño", "CASA", "Roberto Laforcada", "2025-11-09"),
    @("53412356", "Paris Laforcada", "Casa", "Roberto Laforcada", "2025-11-09"),
    @("35115812", "Lisandro Laforcada", "Casa", "Roberto Laforcada", "2025-11-09"),
    @("59610581", "Patrick Laforcada", "Casa", "Roberto Laforcada", "2025-11-09")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = $j + 1
        $cell = $ws.Cells.Item($row, $col)
        # Force text storage (matches source file, which keeps numeric IDs and
        # dates as literal text) instead of letting Excel auto-convert
        # numeric-looking / date-looking strings into numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $data[$i][$j]
        # Reset formatting back to the default (unstyled) cell style so the
        # new data rows don't pick up visible formatting, same as row 2 in
        # the original workbook.
        $cell.Style = "Normal"
    }
}
